# Refresh the crypto symbol list snapshot (Price / Volume(1h) columns).
# Mirrors the scheduled GitHub Actions scrape that updates cryptos.xlsx.
#
# D<row> = Price, E<row> = Volume(1h). Both columns hold text values in the
# source sheet (e.g. "278.39", "1.19%"), not numbers/percentages, so each
# new value is written with a leading apostrophe to force Excel to keep it
# as literal text instead of auto-converting to a number/percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @{ column letter = new value }
$updates = [ordered]@{
    2 = @{ "D"="278.39"; "E"="1.19%" }
    3 = @{ "D"="27.35"; "E"="0.31%" }
    4 = @{ "D"="4.843"; "E"="1.53%" }
    5 = @{ "D"="0.06378"; "E"="0.81%" }
    6 = @{ "D"="7.021"; "E"="1.34%" }
    7 = @{ "D"="1.320"; "E"="-1.99%" }
    8 = @{ "D"="0.8906"; "E"="1.46%" }
    9 = @{ "D"="0.1517"; "E"="0.83%" }
    10 = @{ "D"="0.05492"; "E"="9.25%" }
    11 = @{ "D"="0.07488"; "E"="-0.33%" }
    12 = @{ "D"="0.02962"; "E"="0.62%" }
    13 = @{ "D"="0.08962"; "E"="-0.55%" }
    14 = @{ "D"="0.001581"; "E"="0.73%" }
    15 = @{ "D"="0.0006345"; "E"="-0.14%" }
    16 = @{ "D"="0.006000"; "E"="2.69%" }
    17 = @{ "D"="3.477"; "E"="0.96%" }
    18 = @{ "D"="3.299"; "E"="0.05%" }
    19 = @{ "D"="2.234"; "E"="-1.65%" }
    21 = @{ "D"="0.1348"; "E"="0.43%" }
    22 = @{ "D"="3.910"; "E"="0.01%" }
    23 = @{ "D"="0.1504"; "E"="8.94%" }
    24 = @{ "D"="0.04397"; "E"="-0.37%" }
    25 = @{ "D"="0.001175"; "E"="0.35%" }
    26 = @{ "D"="0.004285"; "E"="11.30%" }
    28 = @{ "E"="-1.83%" }
    29 = @{ "D"="0.0001650"; "E"="-14.71%" }
    40 = @{ "D"="0.04025"; "E"="-2.18%" }
    41 = @{ "D"="0.006683"; "E"="-2.17%" }
    42 = @{ "D"="0.1397"; "E"="18.91%" }
    43 = @{ "D"="0.002065"; "E"="0.81%" }
    44 = @{ "D"="0.01111"; "E"="-3.94%" }
    45 = @{ "D"="0.00005535"; "E"="7.28%" }
    47 = @{ "D"="0.01846"; "E"="-19.68%" }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $newValue = $cols[$col]
        $ws.Range("$col$row").Value = "'$newValue"
    }
}

